$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains a date serial value that was updated
# from 45188 (2023-09-19) to 45189 (2023-09-20) for every data row (2-118).
$ws.Range("C2:C118").Value = 45189
